$wb = $excel.ActiveWorkbook

# --- Controller sheet: mark AttendanceSheetController checklist row (row 2) as
#     done - "Yes" for UnitOfWork DI / Index Order By, "No" for the rest ---
$wsController = $wb.Worksheets.Item("Controller")
[void]$wsController.Activate()
$wsController.Range("B2").Value = "Yes"
$wsController.Range("C2").Value = "Yes"
$wsController.Range("D2:L2").Value = "No"
[void]$wsController.Range("M2").Select()

# --- Views sheet: mark AttendanceSheet checklist row (row 2) as done for the
#     first four checklist columns ---
$wsViews = $wb.Worksheets.Item("Views")
[void]$wsViews.Activate()
$wsViews.Range("B2:E2").Value = "Yes"
[void]$wsViews.Range("F2").Select()

# --- Model sheet: selection moved only ---
$wsModel = $wb.Worksheets.Item("Model")
[void]$wsModel.Activate()
[void]$wsModel.Range("J7").Select()

# Final active sheet is Views (workbook activeTab=2), with its own selection
# left at F2, matching the target state.
[void]$wsViews.Activate()
[void]$wsViews.Range("F2").Select()
